$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3016.2727
$ws.Range("I40").Value = 7389.5
$ws.Range("J40").Value = 2044.4445
$ws.Range("K40").Value = 7389.5
$ws.Range("L40").Value = 2044.4445
$ws.Range("M40").Value = -7214.5
$ws.Range("N40").Value = -2394.4445
$ws.Range("H58").Value = 1669.7255
$ws.Range("I58").Value = 294.14285
$ws.Range("J58").Value = 2190.2163
$ws.Range("K58").Value = 882.4285500000001
$ws.Range("L58").Value = 6570.6489
$ws.Range("M58").Value = -732.4285500000001
$ws.Range("N58").Value = -6870.6489
$ws.Range("H107").Value = 770
$ws.Range("I107").Value = 299.33334
$ws.Range("J107").Value = 1299.5
$ws.Range("K107").Value = 299.33334
$ws.Range("L107").Value = 1299.5
$ws.Range("M107").Value = 1620.66666
$ws.Range("N107").Value = -5139.5
$ws.Range("H111").Value = 5062.5
$ws.Range("I111").Value = 7550
$ws.Range("J111").Value = 2575
$ws.Range("K111").Value = 22650
$ws.Range("L111").Value = 7725
$ws.Range("M111").Value = -19583
$ws.Range("N111").Value = -13859
$ws.Range("H116").Value = 2818.75
$ws.Range("I116").Value = 1712.5
$ws.Range("J116").Value = 3925
$ws.Range("K116").Value = 1712.5
$ws.Range("L116").Value = 3925
$ws.Range("M116").Value = 1729.5
$ws.Range("N116").Value = -10809
$ws.Range("H132").Value = 2383087
$ws.Range("I132").Value = 2859282.2
$ws.Range("J132").Value = 2110.8
$ws.Range("K132").Value = 8577846.600000001
$ws.Range("L132").Value = 6332.400000000001
$ws.Range("M132").Value = -8575316.600000001
$ws.Range("N132").Value = -11392.4
$ws.Range("H137").Value = 1481.1333
$ws.Range("I137").Value = 1131.9131
$ws.Range("J137").Value = 2628.5715
$ws.Range("K137").Value = 3395.7393
$ws.Range("L137").Value = 7885.7145
$ws.Range("M137").Value = -845.7393000000002
$ws.Range("N137").Value = -12985.7145
$ws.Range("H138").Value = 3007.0886
$ws.Range("I138").Value = 2377.3
$ws.Range("J138").Value = 3220.5762
$ws.Range("K138").Value = 7131.900000000001
$ws.Range("L138").Value = 9661.7286
$ws.Range("M138").Value = -1991.900000000001
$ws.Range("N138").Value = -19941.7286

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 33334264
$ws.Range("J2").Value = 1360.2667
$ws.Range("L2").Value = 1360.2667
$ws.Range("N2").Value = -1586.2667
$ws.Range("H61").Value = 1890.7142
$ws.Range("I61").Value = 1548
$ws.Range("K61").Value = 1548
$ws.Range("M61").Value = -1336
$ws.Range("H88").Value = 10162.6
$ws.Range("I88").Value = 20006
$ws.Range("J88").Value = 7701.75
$ws.Range("K88").Value = 20006
$ws.Range("L88").Value = 7701.75
$ws.Range("M88").Value = -19600
$ws.Range("N88").Value = -8513.75
$ws.Range("H91").Value = 10162.6
$ws.Range("I91").Value = 20006
$ws.Range("J91").Value = 7701.75
$ws.Range("K91").Value = 20006
$ws.Range("L91").Value = 7701.75
$ws.Range("M91").Value = -18602
$ws.Range("N91").Value = -10509.75
$ws.Range("H116").Value = 33334264
$ws.Range("J116").Value = 1360.2667
$ws.Range("L116").Value = 1360.2667
$ws.Range("N116").Value = -5948.2667
$ws.Range("H136").Value = 1890.7142
$ws.Range("I136").Value = 1548
$ws.Range("K136").Value = 4644
$ws.Range("M136").Value = -2094

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 33334264
$ws.Range("J3").Value = 1360.2667
$ws.Range("L3").Value = 1360.2667
$ws.Range("N3").Value = -1588.2667
$ws.Range("H105").Value = 2029.909
$ws.Range("I105").Value = 2029.909
$ws.Range("K105").Value = 2029.909
$ws.Range("M105").Value = -282.9090000000001
$ws.Range("H107").Value = 1735.0312
$ws.Range("I107").Value = 1516.9474
$ws.Range("J107").Value = 2053.7693
$ws.Range("K107").Value = 1516.9474
$ws.Range("L107").Value = 2053.7693
$ws.Range("M107").Value = 403.0526
$ws.Range("N107").Value = -5893.7693

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 618.2381
$ws.Range("I107").Value = 590.5854
$ws.Range("K107").Value = 590.5854
$ws.Range("M107").Value = 1329.4146
$ws.Range("H122").Value = 901.3077
$ws.Range("I122").Value = 826
$ws.Range("J122").Value = 1021.8
$ws.Range("K122").Value = 2478
$ws.Range("L122").Value = 3065.4
$ws.Range("M122").Value = -28
$ws.Range("N122").Value = -7965.4

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 2332.7222
$ws.Range("I140").Value = 1440.5294
$ws.Range("J140").Value = 3131
$ws.Range("K140").Value = 4321.5882
$ws.Range("L140").Value = 9393
$ws.Range("M140").Value = 858.4117999999999
$ws.Range("N140").Value = -19753

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 909.875
$ws.Range("I16").Value = 546.5
$ws.Range("K16").Value = 546.5
$ws.Range("M16").Value = -376.5
$ws.Range("H22").Value = 359
$ws.Range("I22").Value = 492.2857
$ws.Range("J22").Value = 225.71428
$ws.Range("K22").Value = 492.2857
$ws.Range("L22").Value = 225.71428
$ws.Range("M22").Value = -197.2857
$ws.Range("N22").Value = -815.71428
$ws.Range("H27").Value = 359
$ws.Range("I27").Value = 492.2857
$ws.Range("J27").Value = 225.71428
$ws.Range("K27").Value = 492.2857
$ws.Range("L27").Value = 225.71428
$ws.Range("M27").Value = -385.2857
$ws.Range("N27").Value = -439.71428
$ws.Range("H61").Value = 23811196
$ws.Range("I61").Value = 1722.2222
$ws.Range("J61").Value = 66668250
$ws.Range("K61").Value = 1722.2222
$ws.Range("L61").Value = 66668250
$ws.Range("M61").Value = -1520.2222
$ws.Range("N61").Value = -66668654
$ws.Range("H113").Value = 23811196
$ws.Range("I113").Value = 1722.2222
$ws.Range("J113").Value = 66668250
$ws.Range("K113").Value = 1722.2222
$ws.Range("L113").Value = 66668250
$ws.Range("M113").Value = 447.7778000000001
$ws.Range("N113").Value = -66672590
$ws.Range("H122").Value = 3571.8215
$ws.Range("I122").Value = 4373.8
$ws.Range("K122").Value = 13121.4
$ws.Range("M122").Value = -10671.4
$ws.Range("H136").Value = 4733.1577
$ws.Range("I136").Value = 5079
$ws.Range("K136").Value = 15237
$ws.Range("M136").Value = -12687

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 41674480
$ws.Range("I126").Value = 62509636
$ws.Range("K126").Value = 187528908
$ws.Range("M126").Value = -187526438
